$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.620.95'
$ws.Range("E2").Value = '  +0.36%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.864.64'
$ws.Range("E3").Value = '  +1.11%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '326.20'
$ws.Range("E5").Value = '  -2.48%  '

$ws.Range("E6").Value = '  +0.03%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4648'
$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3896'
$ws.Range("E8").Value = '  +0.11%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07908'
$ws.Range("E9").Value = '  +0.39%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9740'
$ws.Range("E10").Value = '  +0.06%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '22.29'
$ws.Range("E11").Value = '  +0.90%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.888.47'
$ws.Range("E12").Value = '  +2.34%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.740'
$ws.Range("E13").Value = '  -1.09%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.956'
$ws.Range("E14").Value = '  +0.22%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06919'
$ws.Range("E15").Value = '  -0.52%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '89.06'
$ws.Range("E16").Value = '  +1.41%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.004'
$ws.Range("E17").Value = '  +0.00%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001003'
$ws.Range("E18").Value = '  +0.21%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.88'
$ws.Range("E19").Value = '  -0.68%  '

$ws.Range("E20").Value = '  +0.14%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '28.591.63'
$ws.Range("E21").Value = '  +0.16%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.334'
$ws.Range("E22").Value = '  -0.47%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.10'
$ws.Range("E23").Value = '  -0.50%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.110'
$ws.Range("E24").Value = '  -2.94%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.031.37'
$ws.Range("E25").Value = '  -1.57%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '154.97'
$ws.Range("E26").Value = '  +1.05%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.33'
$ws.Range("E27").Value = '  -0.15%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.819'
$ws.Range("E28").Value = '  -2.25%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.999'
$ws.Range("E29").Value = '  +0.18%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '119.54'
$ws.Range("E30").Value = '  +1.83%  '

$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.9429'
$ws.Range("E31").Value = '  -2.07%  '

$ws.Range("B32").Value = 'Stellar'
$ws.Range("C32").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09322'
$ws.Range("E32").Value = '  -0.58%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.319'
$ws.Range("E33").Value = '  -0.62%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.339'
$ws.Range("E34").Value = '  -0.05%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.334'
$ws.Range("E35").Value = '  -3.50%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05865'
$ws.Range("E36").Value = '  -3.29%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02132'
$ws.Range("E37").Value = '  -2.46%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.149'
$ws.Range("E38").Value = '  -1.14%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '7.837'
$ws.Range("E39").Value = '  +2.74%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5654'
$ws.Range("E40").Value = '  -0.25%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '9.982'
$ws.Range("E41").Value = '  -1.07%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1776'
$ws.Range("E42").Value = '  -0.81%  '

$ws.Range("E43").Value = '  +4.16%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '11.70'
$ws.Range("E44").Value = '  -0.33%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5320'
$ws.Range("E45").Value = '  -0.60%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.170'
$ws.Range("E46").Value = '  -9.26%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.140'
$ws.Range("E47").Value = '  -8.76%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.851'
$ws.Range("E48").Value = '  -1.95%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '114.19'
$ws.Range("E49").Value = '  +0.77%  '

$ws.Range("B50").Value = 'MXToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.346'
$ws.Range("E50").Value = '  +0.63%  '

$ws.Range("B51").Value = 'PaxDollar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.004'
$ws.Range("E51").Value = '  +0.22%  '
